$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Individuals")

# Add a new test individual ("Individual_with_NAs") whose numerical
# characteristics (Weight, Height, Age) are intentionally left blank to
# represent missing/NA values, matching the columns already used by the
# other rows (IndividualId, Species, Population, Gender).
$ws.Range("A4").Value = "Individual_with_NAs"
$ws.Range("B4").Value = "Human"
$ws.Range("C4").Value = "European_ICRP_2002"
$ws.Range("D4").Value = "MALE"

# Leave E4:G4 (Weight/Height/Age) empty to represent the NA values.

# Select the newly populated range, mirroring the active cell/selection
# left behind after entering this row.
$ws.Range("B4:D4").Select()
